$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column B (usageTypeId) before the current "name" column
$ws.Columns("B:B").Insert()

# Header row
$ws.Range("B1").Value = "usageTypeId"

# usageTypeId values for data rows 2-28
$values = @(3,3,3,3,3,3,3,3,3,4,4,4,4,4,4,4,4,4,6,6,6,6,6,6,6,6,6)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("C6").Select()
